# Rewrites Janna Gardner's resume body into a plainer/standard layout.
# The whole resume body lives in ONE paragraph / ONE run, with the visual
# "lines" separated by <w:br/> manual line breaks (not paragraph marks).
# Word's Find/Replace treats a manual line break as the "^l" special
# character, and it can match/insert "^l" across existing <w:t> run
# boundaries, so we do the whole rewrite as a single Find.Execute() call
# built from the old/new "lines" joined with "^l".

$d = $word.ActiveDocument

$oldLines = @(
    "J a n n a G a r d n e r",
    "4567 Main Street, Chicago, Illinois 98052 • (716) 555-0100 • janna@example.com",
    "Software Engineer with 6+ years of experience assisting with and fulfilling organization's technical needs and requirements. A proven track record of using excellent personal, communication, and organizational skills to lead and improve software development teams. Proficient in project management software and data analytics. Possess excellent communication skills, high quality of work, and is driven and highly self-motivated. Team player with the ability to work independently.",
    "",
    "Experience",
    "20XX – PRESENT",
    "Software Engineer | Lamna Healthcare Company | Chicago, Illinois",
    "Review, update, and revise company hiring practices, vacation, and other human resources policies to ensure compliance with OSHA and all local, state, and federal labor regulations. Led the development team to build and deploy a dedicated recruitment website.",
    "",
    "JUNE 20XX – AUGUST 20XX",
    "Intern | Wholeness Healthcare | Boomtown, Ohio",
    "Assisted in recruitment outreach to prospective employees.",
    "",
    "Skills",
    "Proficient with project management software",
    "Data analytics",
    "Excellent time management skills",
    "Conflict management",
    "Team player",
    "",
    "Education",
    "MAY 20XX",
    "Bachelor of Arts Human Resources Management | Jasper University | Ft. Lauderdale, FL",
    "3.8 GPA • Member of university’s Honor Society",
    "",
    "Activities",
    "Literature • Environmental conservation • Art • Yoga • Skiing • Travel"
)

$newLines = @(
    "Jane Gardner",
    "Contact: (716) 555-0100",
    "Email: janna@example.com",
    "",
    "Software Skills: Proficient in various project management software, data analytics tools",
    "Other Skills: Fast typist (96 wpm), excellent time management, public speaking, conflict management",
    "",
    "Education",
    "Bachelor of Arts in Human Resources Management: Jasper University, Ft. Lauderdale, FL, May 20xx - GPA 3.8, Member of University’s Honor Society",
    "",
    "Work Experience",
    "Wholeness Healthcare - Human Resources Intern (June 20xx - August 20xx)",
    "> Assisted in recruitment outreach",
    "> Organized seminars",
    "> Administrative tasks",
    "",
    "Lamna Healthcare Company - Human Resources Generalist (20xx - Present)",
    "> Developed recruitment programs",
    "> Created responsive work environment",
    "> Revise and update company policies",
    "",
    "Hobbies",
    "Literature, Environmental conservation, Art, Yoga, Skiing, and Travel"
)

$findText = [string]::Join("^l", $oldLines)
$replaceText = [string]::Join("^l", $newLines)

$found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)

if (-not $found) {
    throw "Resume body block was not found; nothing was replaced."
}

Write-Output "Resume body rewritten: $found"
